$d = $word.ActiveDocument

$d.Content.Find.Execute("34+23=57", $true, $false, $false, $false, $false, $true, 1, $false, "13+50=63", 2) | Out-Null
$d.Content.Find.Execute("55-35=20", $true, $false, $false, $false, $false, $true, 1, $false, "98-71=27", 2) | Out-Null
$d.Content.Find.Execute("16+31=47", $true, $false, $false, $false, $false, $true, 1, $false, "74-58=16", 2) | Out-Null
$d.Content.Find.Execute("52-4=48", $true, $false, $false, $false, $false, $true, 1, $false, "94-41=53", 2) | Out-Null
$d.Content.Find.Execute("99-34=65", $true, $false, $false, $false, $false, $true, 1, $false, "23+7=30", 2) | Out-Null
$d.Content.Find.Execute("33-1=32", $true, $false, $false, $false, $false, $true, 1, $false, "51-44=7", 2) | Out-Null
$d.Content.Find.Execute("7+21=28", $true, $false, $false, $false, $false, $true, 1, $false, "75+0=75", 2) | Out-Null
$d.Content.Find.Execute("79-35=44", $true, $false, $false, $false, $false, $true, 1, $false, "82-2=80", 2) | Out-Null
$d.Content.Find.Execute("53+30=83", $true, $false, $false, $false, $false, $true, 1, $false, "12+3=15", 2) | Out-Null
$d.Content.Find.Execute("31-19=12", $true, $false, $false, $false, $false, $true, 1, $false, "64-16=48", 2) | Out-Null
$d.Content.Find.Execute("1+35=36", $true, $false, $false, $false, $false, $true, 1, $false, "74-3=71", 2) | Out-Null
$d.Content.Find.Execute("73-17=56", $true, $false, $false, $false, $false, $true, 1, $false, "76-72=4", 2) | Out-Null
$d.Content.Find.Execute("97-72=25", $true, $false, $false, $false, $false, $true, 1, $false, "31+9=40", 2) | Out-Null
$d.Content.Find.Execute("46-22=24", $true, $false, $false, $false, $false, $true, 1, $false, "76-24=52", 2) | Out-Null
$d.Content.Find.Execute("23-6=17", $true, $false, $false, $false, $false, $true, 1, $false, "19+7=26", 2) | Out-Null
$d.Content.Find.Execute("2+35=37", $true, $false, $false, $false, $false, $true, 1, $false, "77-37=40", 2) | Out-Null
$d.Content.Find.Execute("45+40=85", $true, $false, $false, $false, $false, $true, 1, $false, "41+11=52", 2) | Out-Null
$d.Content.Find.Execute("25-13=12", $true, $false, $false, $false, $false, $true, 1, $false, "95-74=21", 2) | Out-Null
$d.Content.Find.Execute("67-25=42", $true, $false, $false, $false, $false, $true, 1, $false, "67-46=21", 2) | Out-Null
$d.Content.Find.Execute("43+33=76", $true, $false, $false, $false, $false, $true, 1, $false, "96-59=37", 2) | Out-Null
$d.Content.Find.Execute("2+40=42", $true, $false, $false, $false, $false, $true, 1, $false, "38+48=86", 2) | Out-Null
$d.Content.Find.Execute("67+9=76", $true, $false, $false, $false, $false, $true, 1, $false, "80+15=95", 2) | Out-Null
$d.Content.Find.Execute("68-34=34", $true, $false, $false, $false, $false, $true, 1, $false, "32+28=60", 2) | Out-Null
$d.Content.Find.Execute("46-16=30", $true, $false, $false, $false, $false, $true, 1, $false, "57+7=64", 2) | Out-Null
$d.Content.Find.Execute("55-6=49", $true, $false, $false, $false, $false, $true, 1, $false, "47-26=21", 2) | Out-Null
$d.Content.Find.Execute("43+25=68", $true, $false, $false, $false, $false, $true, 1, $false, "90-14=76", 2) | Out-Null
$d.Content.Find.Execute("56-9=47", $true, $false, $false, $false, $false, $true, 1, $false, "89-73=16", 2) | Out-Null
$d.Content.Find.Execute("69+27=96", $true, $false, $false, $false, $false, $true, 1, $false, "15+8=23", 2) | Out-Null
$d.Content.Find.Execute("99-7=92", $true, $false, $false, $false, $false, $true, 1, $false, "4+12=16", 2) | Out-Null
$d.Content.Find.Execute("51+23=74", $true, $false, $false, $false, $false, $true, 1, $false, "21-9=12", 2) | Out-Null
$d.Content.Find.Execute("45-4=41", $true, $false, $false, $false, $false, $true, 1, $false, "91-9=82", 2) | Out-Null
$d.Content.Find.Execute("1+46=47", $true, $false, $false, $false, $false, $true, 1, $false, "18+64=82", 2) | Out-Null
$d.Content.Find.Execute("30+33=63", $true, $false, $false, $false, $false, $true, 1, $false, "99-18=81", 2) | Out-Null
$d.Content.Find.Execute("99-73=26", $true, $false, $false, $false, $false, $true, 1, $false, "11+10=21", 2) | Out-Null
$d.Content.Find.Execute("3+60=63", $true, $false, $false, $false, $false, $true, 1, $false, "69+13=82", 2) | Out-Null
$d.Content.Find.Execute("89-19=70", $true, $false, $false, $false, $false, $true, 1, $false, "69+7=76", 2) | Out-Null
$d.Content.Find.Execute("5+43=48", $true, $false, $false, $false, $false, $true, 1, $false, "72+9=81", 2) | Out-Null
$d.Content.Find.Execute("97-47=50", $true, $false, $false, $false, $false, $true, 1, $false, "48-28=20", 2) | Out-Null
$d.Content.Find.Execute("26-23=3", $true, $false, $false, $false, $false, $true, 1, $false, "80-35=45", 2) | Out-Null
$d.Content.Find.Execute("86-22=64", $true, $false, $false, $false, $false, $true, 1, $false, "61+31=92", 2) | Out-Null
$d.Content.Find.Execute("4+54=58", $true, $false, $false, $false, $false, $true, 1, $false, "16-8=8", 2) | Out-Null
$d.Content.Find.Execute("51-5=46", $true, $false, $false, $false, $false, $true, 1, $false, "16+81=97", 2) | Out-Null
$d.Content.Find.Execute("81-81=0", $true, $false, $false, $false, $false, $true, 1, $false, "91-28=63", 2) | Out-Null
$d.Content.Find.Execute("87-21=66", $true, $false, $false, $false, $false, $true, 1, $false, "47-40=7", 2) | Out-Null
$d.Content.Find.Execute("92-70=22", $true, $false, $false, $false, $false, $true, 1, $false, "2+27=29", 2) | Out-Null
$d.Content.Find.Execute("88-27=61", $true, $false, $false, $false, $false, $true, 1, $false, "11-0=11", 2) | Out-Null
$d.Content.Find.Execute("12+38=50", $true, $false, $false, $false, $false, $true, 1, $false, "5+14=19", 2) | Out-Null
$d.Content.Find.Execute("93-55=38", $true, $false, $false, $false, $false, $true, 1, $false, "12+18=30", 2) | Out-Null
$d.Content.Find.Execute("86-20=66", $true, $false, $false, $false, $false, $true, 1, $false, "29+38=67", 2) | Out-Null
$d.Content.Find.Execute("92-15=77", $true, $false, $false, $false, $false, $true, 1, $false, "57-5=52", 2) | Out-Null
$d.Content.Find.Execute("41+35=76", $true, $false, $false, $false, $false, $true, 1, $false, "75-40=35", 2) | Out-Null
$d.Content.Find.Execute("97-42=55", $true, $false, $false, $false, $false, $true, 1, $false, "98-83=15", 2) | Out-Null
$d.Content.Find.Execute("53-16=37", $true, $false, $false, $false, $false, $true, 1, $false, "95-13=82", 2) | Out-Null
$d.Content.Find.Execute("31-8=23", $true, $false, $false, $false, $false, $true, 1, $false, "50-30=20", 2) | Out-Null
$d.Content.Find.Execute("10+82=92", $true, $false, $false, $false, $false, $true, 1, $false, "85-43=42", 2) | Out-Null
$d.Content.Find.Execute("47-23=24", $true, $false, $false, $false, $false, $true, 1, $false, "63-17=46", 2) | Out-Null
$d.Content.Find.Execute("19-15=4", $true, $false, $false, $false, $false, $true, 1, $false, "84-23=61", 2) | Out-Null
$d.Content.Find.Execute("84-56=28", $true, $false, $false, $false, $false, $true, 1, $false, "27+52=79", 2) | Out-Null
$d.Content.Find.Execute("53-9=44", $true, $false, $false, $false, $false, $true, 1, $false, "53+20=73", 2) | Out-Null
$d.Content.Find.Execute("15+42=57", $true, $false, $false, $false, $false, $true, 1, $false, "12+44=56", 2) | Out-Null
$d.Content.Find.Execute("72-13=59", $true, $false, $false, $false, $false, $true, 1, $false, "16+50=66", 2) | Out-Null
$d.Content.Find.Execute("98-1=97", $true, $false, $false, $false, $false, $true, 1, $false, "89-77=12", 2) | Out-Null
$d.Content.Find.Execute("78+17=95", $true, $false, $false, $false, $false, $true, 1, $false, "85-40=45", 2) | Out-Null
$d.Content.Find.Execute("92-79=13", $true, $false, $false, $false, $false, $true, 1, $false, "56-6=50", 2) | Out-Null
$d.Content.Find.Execute("18+70=88", $true, $false, $false, $false, $false, $true, 1, $false, "18+30=48", 2) | Out-Null
$d.Content.Find.Execute("81-25=56", $true, $false, $false, $false, $false, $true, 1, $false, "13+19=32", 2) | Out-Null
$d.Content.Find.Execute("73-20=53", $true, $false, $false, $false, $false, $true, 1, $false, "70-65=5", 2) | Out-Null
$d.Content.Find.Execute("21+77=98", $true, $false, $false, $false, $false, $true, 1, $false, "98-94=4", 2) | Out-Null
$d.Content.Find.Execute("11+36=47", $true, $false, $false, $false, $false, $true, 1, $false, "96-61=35", 2) | Out-Null
$d.Content.Find.Execute("81-47=34", $true, $false, $false, $false, $false, $true, 1, $false, "9+57=66", 2) | Out-Null
$d.Content.Find.Execute("78+18=96", $true, $false, $false, $false, $false, $true, 1, $false, "84-3=81", 2) | Out-Null
$d.Content.Find.Execute("46-9=37", $true, $false, $false, $false, $false, $true, 1, $false, "35+24=59", 2) | Out-Null
$d.Content.Find.Execute("30+42=72", $true, $false, $false, $false, $false, $true, 1, $false, "83+4=87", 2) | Out-Null
$d.Content.Find.Execute("3+72=75", $true, $false, $false, $false, $false, $true, 1, $false, "37-4=33", 2) | Out-Null
$d.Content.Find.Execute("43+46=89", $true, $false, $false, $false, $false, $true, 1, $false, "91-59=32", 2) | Out-Null
$d.Content.Find.Execute("7+56=63", $true, $false, $false, $false, $false, $true, 1, $false, "92-53=39", 2) | Out-Null
$d.Content.Find.Execute("35-7=28", $true, $false, $false, $false, $false, $true, 1, $false, "34-5=29", 2) | Out-Null
$d.Content.Find.Execute("24+67=91", $true, $false, $false, $false, $false, $true, 1, $false, "87-19=68", 2) | Out-Null
$d.Content.Find.Execute("92-89=3", $true, $false, $false, $false, $false, $true, 1, $false, "33+9=42", 2) | Out-Null
$d.Content.Find.Execute("40+50=90", $true, $false, $false, $false, $false, $true, 1, $false, "28-2=26", 2) | Out-Null
$d.Content.Find.Execute("4+77=81", $true, $false, $false, $false, $false, $true, 1, $false, "83+10=93", 2) | Out-Null
$d.Content.Find.Execute("90-89=1", $true, $false, $false, $false, $false, $true, 1, $false, "68+8=76", 2) | Out-Null
$d.Content.Find.Execute("88-73=15", $true, $false, $false, $false, $false, $true, 1, $false, "5+81=86", 2) | Out-Null
$d.Content.Find.Execute("54+18=72", $true, $false, $false, $false, $false, $true, 1, $false, "74+18=92", 2) | Out-Null
$d.Content.Find.Execute("23-19=4", $true, $false, $false, $false, $false, $true, 1, $false, "90-43=47", 2) | Out-Null
$d.Content.Find.Execute("91-48=43", $true, $false, $false, $false, $false, $true, 1, $false, "81-0=81", 2) | Out-Null
$d.Content.Find.Execute("49-24=25", $true, $false, $false, $false, $false, $true, 1, $false, "73-32=41", 2) | Out-Null
$d.Content.Find.Execute("74+12=86", $true, $false, $false, $false, $false, $true, 1, $false, "86-43=43", 2) | Out-Null
$d.Content.Find.Execute("29+15=44", $true, $false, $false, $false, $false, $true, 1, $false, "0+81=81", 2) | Out-Null
$d.Content.Find.Execute("33+37=70", $true, $false, $false, $false, $false, $true, 1, $false, "67-40=27", 2) | Out-Null
$d.Content.Find.Execute("42-3=39", $true, $false, $false, $false, $false, $true, 1, $false, "69-58=11", 2) | Out-Null
$d.Content.Find.Execute("17+56=73", $true, $false, $false, $false, $false, $true, 1, $false, "79-58=21", 2) | Out-Null
$d.Content.Find.Execute("64+27=91", $true, $false, $false, $false, $false, $true, 1, $false, "82-60=22", 2) | Out-Null
$d.Content.Find.Execute("29+17=46", $true, $false, $false, $false, $false, $true, 1, $false, "13+11=24", 2) | Out-Null
$d.Content.Find.Execute("59-49=10", $true, $false, $false, $false, $false, $true, 1, $false, "14+78=92", 2) | Out-Null
$d.Content.Find.Execute("17+50=67", $true, $false, $false, $false, $false, $true, 1, $false, "11+7=18", 2) | Out-Null
$d.Content.Find.Execute("57-48=9", $true, $false, $false, $false, $false, $true, 1, $false, "51-15=36", 2) | Out-Null
$d.Content.Find.Execute("82-56=26", $true, $false, $false, $false, $false, $true, 1, $false, "58-24=34", 2) | Out-Null
$d.Content.Find.Execute("86-50=36", $true, $false, $false, $false, $false, $true, 1, $false, "78-73=5", 2) | Out-Null
$d.Content.Find.Execute("63+21=84", $true, $false, $false, $false, $false, $true, 1, $false, "45-22=23", 2) | Out-Null
